# Insert a new data row before the current row 129 ("Vega Modelo de Temuco" /
# "Cebollín" weekly price series). This shifts the existing rows 129-249
# down to 130-250 (carrying their values/styles with them, as Excel's
# EntireRow Insert does) and grows the used range to A1:R250. We then fill
# the freshly inserted row 129 with its own record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 129:249 down to 130:250, leaving a blank row 129 behind.
$ws.Rows.Item(129).Insert()

# Populate the new row 129 with its data.
$ws.Range("A129").Value2 = 10
$ws.Range("B129").Value2 = "Vega Modelo de Temuco"
$ws.Range("C129").Value2 = "La Araucanía"
$ws.Range("D129").Value2 = 44512
$ws.Range("E129").Value2 = 9
$ws.Range("F129").Value2 = 100112037
$ws.Range("G129").Value2 = "Cebollín"
$ws.Range("H129").Value2 = "Sin especificar"
$ws.Range("I129").Value2 = "Primera"
$ws.Range("J129").Value2 = 30
$ws.Range("K129").Value2 = 8000
$ws.Range("L129").Value2 = 8000
$ws.Range("M129").Value2 = 8000
$ws.Range("N129").Value2 = "`$/docena de paquetes"
$ws.Range("O129").Value2 = "Provincia de Cautín"
$ws.Range("P129").Value2 = 667
$ws.Range("Q129").Value2 = 12
$ws.Range("R129").Value2 = "Hortaliza"
